$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 159
$ws.Range("I9").Value = 157.33333
$ws.Range("K9").Value = 157.33333
$ws.Range("M9").Value = 11.66667000000001
$ws.Range("H51").Value = 3500
$ws.Range("J51").Value = 3500
$ws.Range("L51").Value = 3500
$ws.Range("N51").Value = -4468
$ws.Range("H76").Value = 3477.8333
$ws.Range("I76").Value = 3528.111
$ws.Range("J76").Value = 3427.5557
$ws.Range("K76").Value = 3528.111
$ws.Range("L76").Value = 3427.5557
$ws.Range("M76").Value = -3213.111
$ws.Range("N76").Value = -4057.5557
$ws.Range("H79").Value = 3477.8333
$ws.Range("I79").Value = 3528.111
$ws.Range("J79").Value = 3427.5557
$ws.Range("K79").Value = 3528.111
$ws.Range("L79").Value = 3427.5557
$ws.Range("M79").Value = -2436.111
$ws.Range("N79").Value = -5611.5557
$ws.Range("H80").Value = 943.375
$ws.Range("I80").Value = 1475
$ws.Range("J80").Value = 766.1667
$ws.Range("K80").Value = 4425
$ws.Range("L80").Value = 2298.5001
$ws.Range("M80").Value = -3427
$ws.Range("N80").Value = -4294.5001
$ws.Range("H83").Value = 943.375
$ws.Range("I83").Value = 1475
$ws.Range("J83").Value = 766.1667
$ws.Range("K83").Value = 13275
$ws.Range("L83").Value = 6895.5003
$ws.Range("M83").Value = -8283
$ws.Range("N83").Value = -16879.5003
$ws.Range("H113").Value = 3999.8333
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3999.8333
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3999.8333
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10507.8333
$ws.Range("H137").Value = 1327.9286
$ws.Range("I137").Value = 1299
$ws.Range("J137").Value = 1501.5
$ws.Range("K137").Value = 3897
$ws.Range("L137").Value = 4504.5
$ws.Range("M137").Value = -1347
$ws.Range("N137").Value = -9604.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2412.5
$ws.Range("I45").Value = 2412.5
$ws.Range("K45").Value = 2412.5
$ws.Range("M45").Value = -2035.5
$ws.Range("H61").Value = 2301.7144
$ws.Range("I61").Value = 1822.4
$ws.Range("K61").Value = 1822.4
$ws.Range("M61").Value = -1610.4
$ws.Range("H74").Value = 2792
$ws.Range("I74").Value = 1074.6666
$ws.Range("K74").Value = 1074.6666
$ws.Range("M74").Value = -200.6666
$ws.Range("H77").Value = 2792
$ws.Range("I77").Value = 1074.6666
$ws.Range("K77").Value = 5373.333000000001
$ws.Range("M77").Value = -1005.333000000001
$ws.Range("H122").Value = 2579.8
$ws.Range("I122").Value = 2633.3333
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 7899.999899999999
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -5449.999899999999
$ws.Range("N122").Value = -12398.5
$ws.Range("H132").Value = 2712.1667
$ws.Range("I132").Value = 2426.3125
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 7278.9375
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -4748.9375
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 2301.7144
$ws.Range("I136").Value = 1822.4
$ws.Range("K136").Value = 5467.200000000001
$ws.Range("M136").Value = -2917.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1587.2593
$ws.Range("I20").Value = 1645.4
$ws.Range("J20").Value = 1421.1428
$ws.Range("K20").Value = 1645.4
$ws.Range("L20").Value = 1421.1428
$ws.Range("M20").Value = -1398.4
$ws.Range("N20").Value = -1915.1428
$ws.Range("H86").Value = 3671.2173
$ws.Range("I86").Value = 4025.353
$ws.Range("J86").Value = 2667.8333
$ws.Range("K86").Value = 4025.353
$ws.Range("L86").Value = 2667.8333
$ws.Range("M86").Value = -2902.353
$ws.Range("N86").Value = -4913.8333
$ws.Range("H89").Value = 3671.2173
$ws.Range("I89").Value = 4025.353
$ws.Range("J89").Value = 2667.8333
$ws.Range("K89").Value = 20126.765
$ws.Range("L89").Value = 13339.1665
$ws.Range("M89").Value = -14510.765
$ws.Range("N89").Value = -24571.1665
$ws.Range("H105").Value = 52633510
$ws.Range("I105").Value = 62501820
$ws.Range("K105").Value = 62501820
$ws.Range("M105").Value = -62500073
$ws.Range("H107").Value = 1854.2
$ws.Range("I107").Value = 1447.1
$ws.Range("J107").Value = 2261.3
$ws.Range("K107").Value = 1447.1
$ws.Range("L107").Value = 2261.3
$ws.Range("M107").Value = 472.9000000000001
$ws.Range("N107").Value = -6101.3

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087.7942
$ws.Range("I31").Value = 1059.2
$ws.Range("J31").Value = 2899.842
$ws.Range("K31").Value = 1059.2
$ws.Range("L31").Value = 2899.842
$ws.Range("M31").Value = -764.2
$ws.Range("N31").Value = -3489.842
$ws.Range("H34").Value = 2087.7942
$ws.Range("I34").Value = 1059.2
$ws.Range("J34").Value = 2899.842
$ws.Range("K34").Value = 1059.2
$ws.Range("L34").Value = 2899.842
$ws.Range("M34").Value = -857.2
$ws.Range("N34").Value = -3303.842
$ws.Range("H58").Value = 1235
$ws.Range("I58").Value = 1265.4
$ws.Range("J58").Value = 1083
$ws.Range("K58").Value = 1265.4
$ws.Range("L58").Value = 1083
$ws.Range("M58").Value = -1062.4
$ws.Range("N58").Value = -1489
$ws.Range("H86").Value = 5575615.5
$ws.Range("I86").Value = 13336794
$ws.Range("J86").Value = 31916.428
$ws.Range("K86").Value = 13336794
$ws.Range("L86").Value = 31916.428
$ws.Range("M86").Value = -13335671
$ws.Range("N86").Value = -34162.428
$ws.Range("H89").Value = 5575615.5
$ws.Range("I89").Value = 13336794
$ws.Range("J89").Value = 31916.428
$ws.Range("K89").Value = 66683970
$ws.Range("L89").Value = 159582.14
$ws.Range("M89").Value = -66678354
$ws.Range("N89").Value = -170814.14
$ws.Range("H94").Value = 1473.5
$ws.Range("I94").Value = 1547
$ws.Range("J94").Value = 1400
$ws.Range("K94").Value = 1547
$ws.Range("L94").Value = 1400
$ws.Range("M94").Value = -1096
$ws.Range("N94").Value = -2302
$ws.Range("H107").Value = 600.087
$ws.Range("I107").Value = 374.0909
$ws.Range("K107").Value = 374.0909
$ws.Range("M107").Value = 1545.9091
$ws.Range("H122").Value = 881.25
$ws.Range("I122").Value = 792.8570999999999
$ws.Range("K122").Value = 2378.5713
$ws.Range("M122").Value = 71.42870000000039
$ws.Range("H136").Value = 1235
$ws.Range("I136").Value = 1265.4
$ws.Range("J136").Value = 1083
$ws.Range("K136").Value = 3796.2
$ws.Range("L136").Value = 3249
$ws.Range("M136").Value = -1246.2
$ws.Range("N136").Value = -8349

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39.666668
$ws.Range("J12").Value = 35
$ws.Range("L12").Value = 105
$ws.Range("N12").Value = -451
$ws.Range("H122").Value = 614.7143
$ws.Range("I122").Value = 449.33334
$ws.Range("K122").Value = 4044.00006
$ws.Range("M122").Value = -1594.00006
$ws.Range("H131").Value = 33334888
$ws.Range("I131").Value = 125000210
$ws.Range("J131").Value = 2044.4546
$ws.Range("K131").Value = 375000630
$ws.Range("L131").Value = 6133.3638
$ws.Range("M131").Value = -374995590
$ws.Range("N131").Value = -16213.3638
$ws.Range("H137").Value = 2335.4167
$ws.Range("I137").Value = 1153.1666
$ws.Range("J137").Value = 3517.6667
$ws.Range("K137").Value = 3459.4998
$ws.Range("L137").Value = 10553.0001
$ws.Range("M137").Value = 1640.5002
$ws.Range("N137").Value = -20753.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18005870
$ws.Range("J70").Value = 18188502
$ws.Range("L70").Value = 18188502
$ws.Range("N70").Value = -18189042
$ws.Range("H73").Value = 18005870
$ws.Range("J73").Value = 18188502
$ws.Range("L73").Value = 18188502
$ws.Range("N73").Value = -18190374
$ws.Range("H102").Value = 4502.4
$ws.Range("I102").Value = 5278
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 5278
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = -3656
$ws.Range("N102").Value = -4644
$ws.Range("H107").Value = 927.0769
$ws.Range("I107").Value = 986
$ws.Range("K107").Value = 986
$ws.Range("M107").Value = 934
$ws.Range("H113").Value = 3300.0715
$ws.Range("I113").Value = 1337
$ws.Range("J113").Value = 3835.4546
$ws.Range("K113").Value = 1337
$ws.Range("L113").Value = 3835.4546
$ws.Range("M113").Value = 833
$ws.Range("N113").Value = -8175.4546
$ws.Range("H122").Value = 252283.33
$ws.Range("J122").Value = 377125
$ws.Range("L122").Value = 1131375
$ws.Range("N122").Value = -1136275

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4079.8
$ws.Range("I40").Value = 3949.5
$ws.Range("K40").Value = 3949.5
$ws.Range("M40").Value = -3813.5
$ws.Range("H122").Value = 25003144
$ws.Range("I122").Value = 62503470
$ws.Range("J122").Value = 2924.5
$ws.Range("K122").Value = 187510410
$ws.Range("L122").Value = 8773.5
$ws.Range("M122").Value = -187507960
$ws.Range("N122").Value = -13673.5
$ws.Range("H123").Value = 19000
$ws.Range("J123").Value = 19000
$ws.Range("L123").Value = 19000
$ws.Range("N123").Value = -28800
